$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label and the value that used to hold a date into an ID/code.
$ws.Range("A1").Value = "编号"
$ws.Range("A2").Value = "'001"

# Move the active selection to A2 (it was F5 before).
$ws.Range("A2").Select()
